# Incorporate updated data from upstream processes through 2024
#
# Sheet1 column E holds the "Solar" series (B1:G1 header row = Biofuel,
# Energy Storage, Natural Gas/Propane, Solar, Waste Gas, Wind) with one
# row per "Open year" (column A). The chart on the sheet plots this same
# range, so updating the cells is the single source of truth the chart
# reads from.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 17 -> Open year 2015, Solar capacity: 11.21 -> 10.37
$ws.Range("E17").Value = 10.37

# Row 26 -> Open year 2024, Solar capacity: 60.84 -> 101.668
$ws.Range("E26").Value = 101.668
